# The commit swaps the raw contents of ppt/theme/theme1.xml (the
# "Office Theme" / "Office" colour scheme, otherwise only wired to the
# Notes Master) and ppt/theme/theme2.xml (the "Integral" / "Red Violet"
# colour scheme that the Slide Master + the deck actually use).  The
# <a:fontScheme> and <a:fmtScheme> blocks in both parts are byte-for-byte
# identical, so the only real content delta is the 12 colours (plus the
# cosmetic theme/colour-scheme names) in <a:clrScheme>.
#
# This host's PowerPoint object model doesn't expose a working
# file-based "load a theme" primitive (Presentation/Master.ApplyTheme,
# Application.OpenThemeFile, ThemeColorScheme.Load/Save are all present
# for API-compat but are no-ops here — there's no real filesystem to
# pull a .thmx from). The supported, live-wired way to edit a theme's
# colours is per swatch, through the Slide's ThemeColorScheme (which
# reads/writes the *Slide Master's* theme part, i.e. theme2.xml; every
# slide shares the one master so any slide's ThemeColorScheme reaches
# the same part):
#
#     Slide.ThemeColorScheme.Colors(i).RGB
#
# with i = 1..12 walking dk1, lt1, dk2, lt2, accent1..accent6, hlink,
# folHlink in that order (verified empirically against this file).
#
# So: push the "Office Theme" palette (currently living in theme1.xml)
# into the live theme (theme2.xml) one swatch at a time.

$p = $ppt.ActivePresentation

# Target palette == the current ppt/theme/theme1.xml <a:clrScheme>
# ("Office"), in Colors(1..12) order.
$officeHex = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

# Any slide's ThemeColorScheme resolves to the single Slide Master's
# theme part, so slide 1 is as good as any.
$tcs = $p.Slides.Item(1).ThemeColorScheme

for ($i = 1; $i -le $officeHex.Count; $i++) {
    $hex = $officeHex[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # COM RGB() packing: 0x00BBGGRR
    $comRgb = $r + ($g * 256) + ($b * 65536)
    $tcs.Colors($i).RGB = $comRgb
}
